# Update the cryptos price/volume table with refreshed values from the source feed.
# D2/D3/... look like plain decimals/money amounts and, if several of them would be
# auto-parsed as numbers by Excel, we force a leading apostrophe so the exact original
# text (e.g. trailing zeros, "1.00") is preserved as a text string, matching the sheet's
# existing inline-string formatting for that column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.837.09"
$ws.Range("E2").Value = "  +3.44%  "
$ws.Range("D3").Value = "3.688.04"
$ws.Range("E3").Value = "  +8.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'588.58"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'179.03"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "3.684.09"
$ws.Range("D8").Value = "'0.623"
$ws.Range("E8").Value = "  +4.85%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.203"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("D11").Value = "'0.615"
$ws.Range("E11").Value = "  +4.73%  "
$ws.Range("D12").Value = "'50.01"
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "4.284.94"
$ws.Range("E14").Value = "  +8.93%  "
$ws.Range("D15").Value = "'683.86"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'8.97"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "71.883.03"
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").Value = "3.671.08"
$ws.Range("E18").Value = "  +8.48%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "'18.20"
$ws.Range("E20").Value = "  +2.79%  "
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "'0.940"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").Value = "'6.16"
$ws.Range("E23").Value = "  +14.77%  "
$ws.Range("D24").Value = "'17.78"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "'103.56"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("E27").Value = "  +5.66%  "
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("D29").Value = "'35.49"
$ws.Range("E29").Value = "  +5.99%  "
$ws.Range("D30").Value = "'9.21"
$ws.Range("E30").Value = "  +5.37%  "
$ws.Range("D31").Value = "'7.38"
$ws.Range("E31").Value = "  +6.10%  "
$ws.Range("D32").Value = "'4.20"
$ws.Range("E32").Value = "  +8.89%  "
$ws.Range("D33").Value = "'574.39"
$ws.Range("E33").Value = "  +4.31%  "
$ws.Range("D34").Value = "'11.32"
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("D36").Value = "'59.81"
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").Value = "3.748.79"
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +3.10%  "
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").Value = "'35.49"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'3.48"
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("E43").Value = "  +8.85%  "
$ws.Range("D44").Value = "'2.79"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("D45").Value = "'0.348"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.88"
$ws.Range("E46").Value = "  +7.72%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.37"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  +3.99%  "
$ws.Range("D49").Value = "'1.43"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "'134.03"
$ws.Range("E51").Value = "  +3.28%  "
